$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "34.481.13"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.809.28"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'225.81"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'0.601"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'36.34"
$ws.Range("E8").Value = "  +3.93%  "
$ws.Range("D9").Value = "'0.294"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'0.0682"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.0964"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "2.066.74"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "'11.37"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "1.805.19"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "34.437.00"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'68.65"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'243.15"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "0.0₃0775"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'11.28"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("E24").Value = "  +5.12%  "
$ws.Range("D25").Value = "'171.33"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'7.90"
$ws.Range("E26").Value = "  +3.83%  "
$ws.Range("D27").Value = "'17.34"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'3.82"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "1.362.18"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "'0.655"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'2.38"
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("D39").Value = "'0.0187"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "'81.24"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").Value = "'2.79"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "'0.939"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "'1.16"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").Value = "'13.43"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'0.0499"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "1.970.46"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'5.83"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'102.66"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'48.46"
$ws.Range("E51").Value = "  -2.21%  "
